# Updated cryptos list on Sat Oct 14 02:47:09 UTC 2023 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns for every coin row (2-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @{ D = <new price text>; E = <new volume text> } (D omitted when unchanged)
$updates = @(
    @{ Row=2; D="26.943.28"; E="  +0.18%  " },
    @{ Row=3; D="1.556.78"; E="  +0.66%  " },
    @{ Row=4; E="  -0.14%  " },
    @{ Row=5; D="206.93"; E="  +0.14%  " },
    @{ Row=6; E="  -0.03%  " },
    @{ Row=7; E="  -0.21%  " },
    @{ Row=8; D="22.06"; E="  +3.01%  " },
    @{ Row=9; E="  +0.31%  " },
    @{ Row=10; D="0.0588"; E="  +0.98%  " },
    @{ Row=11; D="0.0858"; E="  +0.16%  " },
    @{ Row=12; D="1.778.86"; E="  +0.67%  " },
    @{ Row=13; D="1.557.00"; E="  +0.45%  " },
    @{ Row=14; E="  +1.36%  " },
    @{ Row=15; E="  +1.63%  " },
    @{ Row=16; D="26.950.91"; E="  +0.20%  " },
    @{ Row=17; D="61.79"; E="  +0.64%  " },
    @{ Row=18; D="217.79"; E="  +1.61%  " },
    @{ Row=19; E="  +1.89%  " },
    @{ Row=20; D="7.31"; E="  +1.21%  " },
    @{ Row=21; E="  -0.15%  " },
    @{ Row=22; E="  +1.23%  " },
    @{ Row=23; D="9.19"; E="  +0.33%  " },
    @{ Row=24; E="  +0.89%  " },
    @{ Row=25; D="153.71"; E="  +1.30%  " },
    @{ Row=26; D="6.65"; E="  +0.32%  " },
    @{ Row=27; E="  +0.53%  " },
    @{ Row=28; E="  +0.55%  " },
    @{ Row=29; E="  -0.12%  " },
    @{ Row=30; E="  +2.34%  " },
    @{ Row=31; D="1.09"; E="  -0.89%  " },
    @{ Row=32; E="  +0.00%  " },
    @{ Row=33; D="1.422.84"; E="  +4.21%  " },
    @{ Row=34; E="  +4.31%  " },
    @{ Row=35; E="  +3.48%  " },
    @{ Row=36; D="0.980"; E="  +2.16%  " },
    @{ Row=37; E="  +0.08%  " },
    @{ Row=38; E="  -0.13%  " },
    @{ Row=39; E="  +0.02%  " },
    @{ Row=40; D="0.812"; E="  +0.96%  " },
    @{ Row=41; E="  -0.16%  " },
    @{ Row=42; E="  +1.53%  " },
    @{ Row=43; E="  +3.56%  " },
    @{ Row=44; D="0.984"; E="  -0.45%  " },
    @{ Row=45; D="64.71"; E="  +1.96%  " },
    @{ Row=46; E="  +1.37%  " },
    @{ Row=47; D="1.693.04"; E="  +0.66%  " },
    @{ Row=48; D="87.75"; E="  +2.59%  " },
    @{ Row=49; E="  +2.13%  " },
    @{ Row=50; D="0.0$([char]0x2087)0998"; E="  +2.66%  " },
    @{ Row=51; D="0.0959"; E="  +1.28%  " }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("D")) {
        # Leading apostrophe forces literal text (these price strings use
        # "." as a thousands separator and would otherwise be auto-coerced
        # into a number, losing the exact original formatting).
        $ws.Range("D" + $u.Row).Value = "'" + $u.D
        # Excel stamps a quote-prefix style on text-forced cells; put the
        # style back to the workbook default so formatting is unaffected.
        $ws.Range("D" + $u.Row).Style = "Normal"
    }
    if ($u.ContainsKey("E")) {
        $ws.Range("E" + $u.Row).Value = $u.E
    }
}
